# problemaMenorCaminho.xlsx - "resolucao de exercicios sobre o menor caminho"
#
# The sheet models a shortest-path LP (Solver "transportation"-style network
# flow). D2:D8 are the edge-usage decision variables, G2:G7 are the
# flow-balance constraints (SUMIF over inflow/outflow) and D10 is the
# total-cost objective (SUMPRODUCT of edge cost * edge usage).
#
# This edit records a completed solve: the chosen path uses edges 2, 4 and 7
# (D2=D4=D7=1, rest 0), the flow-balance formulas in column G are fixed to
# reference the node id in column F (instead of self-referencing column G,
# which made them vol710atile/circular), and the selection + solver_neg
# defined name are updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Mark which edges are used in the shortest path (decision variables) ---
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 0

# --- Fix the flow-balance formulas in column G: they should look up the
#     node id in column F, not self-reference column G (drops the old
#     always-recalc/circular behaviour). G3:G7 stays one shared formula. ---
$ws.Range("G2").Formula = "=SUMIF(`$B`$2:`$B`$8,F2,`$D`$2:`$D`$8)-SUMIF(`$A`$2:`$A`$8,F2,`$D`$2:`$D`$8)"
$ws.Range("G3:G7").Formula = "=SUMIF(`$B`$2:`$B`$8,F3,`$D`$2:`$D`$8)-SUMIF(`$A`$2:`$A`$8,F3,`$D`$2:`$D`$8)"

# --- Update the solver's "assume non-negative" option (solver_neg 1 -> 2) ---
$wb.Names.Item("solver_neg").RefersToR1C1 = "2"

# --- Match the selection left behind in the saved file ---
$ws.Range("G2:G7").Select()

$wb.Save()
